$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 8266.606
$ws.Cells.Item(51, 10).Value = 8310.299999999999
$ws.Cells.Item(51, 12).Value = 8310.299999999999
$ws.Cells.Item(51, 14).Value = -9278.299999999999

$ws.Cells.Item(112, 8).Value = 1585767.4
$ws.Cells.Item(112, 9).Value = 1250
$ws.Cells.Item(112, 10).Value = 1661220.6
$ws.Cells.Item(112, 11).Value = 3750
$ws.Cells.Item(112, 12).Value = 4983661.800000001
$ws.Cells.Item(112, 13).Value = -2642
$ws.Cells.Item(112, 14).Value = -4985877.800000001

$ws.Cells.Item(121, 8).Value = 5644
$ws.Cells.Item(121, 10).Value = 5644
$ws.Cells.Item(121, 12).Value = 16932
$ws.Cells.Item(121, 14).Value = -20426

$ws.Cells.Item(129, 8).Value = 1411.2858
$ws.Cells.Item(129, 9).Value = 1072.25
$ws.Cells.Item(129, 10).Value = 1863.3334
$ws.Cells.Item(129, 11).Value = 3216.75
$ws.Cells.Item(129, 12).Value = 5590.0002
$ws.Cells.Item(129, 13).Value = 1783.25
$ws.Cells.Item(129, 14).Value = -15590.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2009.4648
$ws.Cells.Item(32, 9).Value = 2060.7761
$ws.Cells.Item(32, 11).Value = 2060.7761
$ws.Cells.Item(32, 13).Value = -1773.7761

$ws.Cells.Item(61, 8).Value = 8865.299999999999
$ws.Cells.Item(61, 9).Value = 3776.5
$ws.Cells.Item(61, 10).Value = 10137.5
$ws.Cells.Item(61, 11).Value = 3776.5
$ws.Cells.Item(61, 12).Value = 10137.5
$ws.Cells.Item(61, 13).Value = -3564.5
$ws.Cells.Item(61, 14).Value = -10561.5

$ws.Cells.Item(74, 8).Value = 78042.5
$ws.Cells.Item(74, 9).Value = 87582.95
$ws.Cells.Item(74, 11).Value = 87582.95
$ws.Cells.Item(74, 13).Value = -86708.95

$ws.Cells.Item(77, 8).Value = 78042.5
$ws.Cells.Item(77, 9).Value = 87582.95
$ws.Cells.Item(77, 11).Value = 437914.75
$ws.Cells.Item(77, 13).Value = -433546.75

$ws.Cells.Item(136, 8).Value = 8865.299999999999
$ws.Cells.Item(136, 9).Value = 3776.5
$ws.Cells.Item(136, 10).Value = 10137.5
$ws.Cells.Item(136, 11).Value = 11329.5
$ws.Cells.Item(136, 12).Value = 30412.5
$ws.Cells.Item(136, 13).Value = -8779.5
$ws.Cells.Item(136, 14).Value = -35512.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 4515.1904
$ws.Cells.Item(99, 9).Value = 4014
$ws.Cells.Item(99, 11).Value = 4014
$ws.Cells.Item(99, 13).Value = -2516

$ws.Cells.Item(105, 8).Value = 2789.4167
$ws.Cells.Item(105, 9).Value = 2588.4546
$ws.Cells.Item(105, 11).Value = 2588.4546
$ws.Cells.Item(105, 13).Value = -841.4546

$ws.Cells.Item(134, 8).Value = 3819.5
$ws.Cells.Item(134, 9).Value = 3939.5
$ws.Cells.Item(134, 10).Value = 2499.5
$ws.Cells.Item(134, 11).Value = 11818.5
$ws.Cells.Item(134, 12).Value = 7498.5
$ws.Cells.Item(134, 13).Value = -9283.5
$ws.Cells.Item(134, 14).Value = -12568.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1645.2727
$ws.Cells.Item(58, 9).Value = 1698
$ws.Cells.Item(58, 11).Value = 1698
$ws.Cells.Item(58, 13).Value = -1495

$ws.Cells.Item(62, 8).Value = 3804.4285
$ws.Cells.Item(62, 9).Value = 2513.7273
$ws.Cells.Item(62, 11).Value = 2513.7273
$ws.Cells.Item(62, 13).Value = -1889.7273

$ws.Cells.Item(65, 8).Value = 3804.4285
$ws.Cells.Item(65, 9).Value = 2513.7273
$ws.Cells.Item(65, 11).Value = 12568.6365
$ws.Cells.Item(65, 13).Value = -9448.636500000001

$ws.Cells.Item(68, 8).Value = 41343.75
$ws.Cells.Item(68, 10).Value = 41343.75
$ws.Cells.Item(68, 12).Value = 41343.75
$ws.Cells.Item(68, 14).Value = -42841.75

$ws.Cells.Item(71, 8).Value = 41343.75
$ws.Cells.Item(71, 10).Value = 41343.75
$ws.Cells.Item(71, 12).Value = 124031.25
$ws.Cells.Item(71, 14).Value = -131519.25

$ws.Cells.Item(94, 8).Value = 1179.5883
$ws.Cells.Item(94, 9).Value = 888
$ws.Cells.Item(94, 11).Value = 888
$ws.Cells.Item(94, 13).Value = -437

$ws.Cells.Item(99, 8).Value = 291193.84
$ws.Cells.Item(99, 9).Value = 504436.16
$ws.Cells.Item(99, 10).Value = 6870.7334
$ws.Cells.Item(99, 11).Value = 504436.16
$ws.Cells.Item(99, 12).Value = 6870.7334
$ws.Cells.Item(99, 13).Value = -502938.16
$ws.Cells.Item(99, 14).Value = -9866.733400000001

$ws.Cells.Item(126, 8).Value = 291193.84
$ws.Cells.Item(126, 9).Value = 504436.16
$ws.Cells.Item(126, 10).Value = 6870.7334
$ws.Cells.Item(126, 11).Value = 1513308.48
$ws.Cells.Item(126, 12).Value = 20612.2002
$ws.Cells.Item(126, 13).Value = -1510838.48
$ws.Cells.Item(126, 14).Value = -25552.2002

$ws.Cells.Item(132, 8).Value = 4880.913
$ws.Cells.Item(132, 9).Value = 4358.6055
$ws.Cells.Item(132, 11).Value = 13075.8165
$ws.Cells.Item(132, 13).Value = -10545.8165

$ws.Cells.Item(134, 8).Value = 1976.1562
$ws.Cells.Item(134, 9).Value = 1732.3448
$ws.Cells.Item(134, 10).Value = 4333
$ws.Cells.Item(134, 11).Value = 5197.0344
$ws.Cells.Item(134, 12).Value = 12999
$ws.Cells.Item(134, 13).Value = -2662.0344
$ws.Cells.Item(134, 14).Value = -18069

$ws.Cells.Item(136, 8).Value = 1645.2727
$ws.Cells.Item(136, 9).Value = 1698
$ws.Cells.Item(136, 11).Value = 5094
$ws.Cells.Item(136, 13).Value = -2544

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1198.129
$ws.Cells.Item(34, 9).Value = 75.333336
$ws.Cells.Item(34, 11).Value = 226.000008
$ws.Cells.Item(34, 13).Value = -142.000008

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 14738
$ws.Cells.Item(43, 9).Value = 3783.2
$ws.Cells.Item(43, 11).Value = 3783.2
$ws.Cells.Item(43, 13).Value = -3632.2

$ws.Cells.Item(46, 8).Value = 18981.6
$ws.Cells.Item(46, 10).Value = 69999
$ws.Cells.Item(46, 12).Value = 69999
$ws.Cells.Item(46, 14).Value = -70311

$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 626.8182
$ws.Cells.Item(22, 9).Value = 579.8570999999999
$ws.Cells.Item(22, 10).Value = 709
$ws.Cells.Item(22, 11).Value = 579.8570999999999
$ws.Cells.Item(22, 12).Value = 709
$ws.Cells.Item(22, 13).Value = -284.8570999999999
$ws.Cells.Item(22, 14).Value = -1299

$ws.Cells.Item(27, 8).Value = 626.8182
$ws.Cells.Item(27, 9).Value = 579.8570999999999
$ws.Cells.Item(27, 10).Value = 709
$ws.Cells.Item(27, 11).Value = 579.8570999999999
$ws.Cells.Item(27, 12).Value = 709
$ws.Cells.Item(27, 13).Value = -472.8570999999999
$ws.Cells.Item(27, 14).Value = -923

$ws.Cells.Item(40, 8).Value = 2103.0476
$ws.Cells.Item(40, 9).Value = 2118.45
$ws.Cells.Item(40, 11).Value = 2118.45
$ws.Cells.Item(40, 13).Value = -1982.45

$ws.Cells.Item(100, 8).Value = 13739.4
$ws.Cells.Item(100, 9).Value = 13424.375
$ws.Cells.Item(100, 11).Value = 13424.375
$ws.Cells.Item(100, 13).Value = -12883.375

$ws.Cells.Item(107, 8).Value = 502770
$ws.Cells.Item(107, 9).Value = 502770
$ws.Cells.Item(107, 11).Value = 502770
$ws.Cells.Item(107, 13).Value = -500850

$ws.Cells.Item(132, 8).Value = 3730.8333
$ws.Cells.Item(132, 9).Value = 3728.0435
$ws.Cells.Item(132, 10).Value = 3795
$ws.Cells.Item(132, 11).Value = 11184.1305
$ws.Cells.Item(132, 12).Value = 11385
$ws.Cells.Item(132, 13).Value = -8654.130500000001
$ws.Cells.Item(132, 14).Value = -16445

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 14).ClearContents()

$ws.Cells.Item(49, 8).Value = 32222.223
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 3488.25
$ws.Cells.Item(122, 9).Value = 3488.25
$ws.Cells.Item(122, 11).Value = 10464.75
$ws.Cells.Item(122, 13).Value = -8014.75

$ws.Cells.Item(132, 8).Value = 2190.3513
$ws.Cells.Item(132, 9).Value = 1851.697
$ws.Cells.Item(132, 10).Value = 2462.9268
$ws.Cells.Item(132, 11).Value = 5555.090999999999
$ws.Cells.Item(132, 12).Value = 7388.780400000001
$ws.Cells.Item(132, 13).Value = -3025.090999999999
$ws.Cells.Item(132, 14).Value = -12448.7804

$ws.Cells.Item(136, 8).Value = 557997.75
$ws.Cells.Item(136, 9).Value = 627247.5600000001
$ws.Cells.Item(136, 10).Value = 3999.5
$ws.Cells.Item(136, 11).Value = 1881742.68
$ws.Cells.Item(136, 12).Value = 11998.5
$ws.Cells.Item(136, 13).Value = -1879192.68
$ws.Cells.Item(136, 14).Value = -17098.5
